$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure the new DataDate values are stored as text (matching the
# existing "YYYYMMDD" text entries in column A), not auto-converted to numbers.
$ws.Range("A32:A36").NumberFormat = "@"

# Append new rows for DataDate "20210120" following the same pattern used
# for the other dates (one row per broker alias: hait, huat, swhy, gtja, zx)
$ws.Range("A32").Value = "20210120"
$ws.Range("B32").Value = "hait"
$ws.Range("C32").Value = "DataFileNotExists"

$ws.Range("A33").Value = "20210120"
$ws.Range("B33").Value = "huat"
$ws.Range("C33").Value = "DataFileNotExists"

$ws.Range("A34").Value = "20210120"
$ws.Range("B34").Value = "swhy"
$ws.Range("C34").Value = "DataFileNotExists"

$ws.Range("A35").Value = "20210120"
$ws.Range("B35").Value = "gtja"
$ws.Range("C35").Value = 23

$ws.Range("A36").Value = "20210120"
$ws.Range("B36").Value = "zx"
$ws.Range("C36").Value = 38

# Restore the plain "Normal" cell style (no explicit style index) on the new
# DataDate cells, matching the unstyled cells used throughout the data rows.
$ws.Range("A32:A36").Style = "Normal"
